$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9806808829307556
$ws.Range("B1").Value = 1.31063187122345
$ws.Range("C1").Value = 2.118067502975464
$ws.Range("D1").Value = 4.560078144073486
$ws.Range("E1").Value = 2.11015248298645
